$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue {
    param($Address, $Text)
    $r = $ws.Range($Address)
    $r.NumberFormat = "@"
    $r.Value = $Text
    $r.Style = "Normal"
}

# ---- Row 2 ----
$ws.Range("A2").Value = "(547, 362)"
$ws.Range("B2").Value = "(484, 381)"
Set-TextValue "C2" "65.80273550544841"
$ws.Range("D2").Value = "0:00:00.690796"
$ws.Range("E2").Value = "[[505, 373, datetime.timedelta(microseconds=136001), 200.84810798557422], [492, 371, datetime.timedelta(microseconds=277003), 21.796434777706118, -646.3889315562218], [486, 385, datetime.timedelta(microseconds=420800), 16.61556618107878, -12.311950087042153], [484, 381, datetime.timedelta(microseconds=565304), 3.631447940954, -22.968382038911415]]"
$ws.Range("F2").Value = "2022-07-18 13:59:41.543855"
Set-TextValue "G2" "30.205834009127525"
$ws.Range("H2").Value = "miss"
$ws.Range("I2").Value = "(428.0, 360.0)"
Set-TextValue "J2" "59.80802621722272"

# ---- Row 3 ----
$ws.Range("A3").Value = "(710, 375)"
$ws.Range("B3").Value = "(844, 415)"
Set-TextValue "C3" "139.84276885130672"
$ws.Range("D3").Value = "0:00:00.590395"
$ws.Range("E3").Value = "[[796, 406, datetime.timedelta(microseconds=208394), 537.4530806935257], [842, 415, datetime.timedelta(microseconds=420094), 51.217152911941916, -1157.4455426204227]]"
$ws.Range("F3").Value = "2022-07-18 13:59:48.176819"
Set-TextValue "G3" "64.19288546066008"
$ws.Range("H3").Value = "hit"
$ws.Range("I3").Value = "(852.0, 360.0)"
Set-TextValue "J3" "55.57877292636101"

# ---- Row 4 ----
$ws.Range("A4").Value = "(559, 348)"
$ws.Range("B4").Value = "(678, 227)"
Set-TextValue "C4" "169.7115199389835"
$ws.Range("D4").Value = "0:00:03.796694"
$ws.Range("E4").Value = "[[477, 323, datetime.timedelta(microseconds=161609), 243.4980367780814], [381, 318, datetime.timedelta(microseconds=330717), 133.4288824928215, -332.8197651927779], [362, 321, datetime.timedelta(microseconds=488226), 18.085346825348815, -236.25029324016475], [362, 335, datetime.timedelta(microseconds=637561), 10.079829261978675, -12.55647312707355], [360, 319, datetime.timedelta(microseconds=792517), 9.339528901905195, -0.9341129087117122], [362, 327, datetime.timedelta(microseconds=934024), 4.052689245828128, -5.660282451068782], [360, 321, datetime.timedelta(seconds=1, microseconds=77222), 37.5954972730523, 434.3685481757035], [362, 335, datetime.timedelta(seconds=1, microseconds=222800), 29.137124833128706, -37.963969658543974], [509, 327, datetime.timedelta(seconds=1, microseconds=373520), 180.92248235523147, 406.36473956442165], [691, 381, datetime.timedelta(seconds=1, microseconds=662303), 131.5777792603786, -74.50472532187362], [465, 325, datetime.timedelta(seconds=1, microseconds=806844), 132.46618361503576, 1.101085655538313], [318, 275, datetime.timedelta(seconds=1, microseconds=956519), 74.51485842321786, -60.5856498321705], [459, 270, datetime.timedelta(seconds=2, microseconds=113482), 570.7052950603302, 4372.415331392753], [772, 319, datetime.timedelta(seconds=2, microseconds=271750), 535.154636715102, -130.8211898628451], [632, 291, datetime.timedelta(seconds=2, microseconds=419753), 156.1341056339755, -902.9608628910967], [296, 241, datetime.timedelta(seconds=2, microseconds=609953), 255.6500541923768, 163.1534701172079], [413, 132, datetime.timedelta(seconds=2, microseconds=775720), 94.62529763159239, -207.5810299602749], [683, 199, datetime.timedelta(seconds=2, microseconds=928935), 137.4678607749794, 46.12008713568444], [670, 212, datetime.timedelta(seconds=3, microseconds=80827), 104.4116054866245, -408.9754078260349], [658, 233, datetime.timedelta(seconds=3, microseconds=221791), 50.058853335468065, -245.06292929449992], [683, 239, datetime.timedelta(seconds=3, microseconds=363827), 32.437896790424865, -48.43223989710275], [672, 224, datetime.timedelta(seconds=3, microseconds=504925), 16.910562681868615, -30.751763348133387], [676, 224, datetime.timedelta(seconds=3, microseconds=654045), 2.8073673498203537, -21.563035161263006]]"
$ws.Range("F4").Value = "2022-07-18 14:00:03.569802"
Set-TextValue "G4" "77.90372180331653"
$ws.Range("H4").Value = "miss"
$ws.Range("I4").Value = "(428.0, 360.0)"
Set-TextValue "J4" "283.17662332897464"

# ---- Row 5 ----
$ws.Range("A5").Value = "Subject Code:"
Set-TextValue "B5" "4"
$ws.Range("C5:J5").ClearContents()

# ---- Row 6: remove entirely ----
$ws.Range("A6:B6").EntireRow.Delete()
